$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells before writing, so that
# numeric-looking values (e.g. "0.999", "7.32") are stored as text,
# matching the source data which uses inline string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = '68.057.55'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '3.550.02'
$ws.Range("E3").Value = '  +1.68%  '
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '618.57'
$ws.Range("E5").Value = '  +1.93%  '
$ws.Range("D6").Value = '154.58'
$ws.Range("E6").Value = '  +3.89%  '
$ws.Range("D7").Value = '3.547.99'
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = '0.489'
$ws.Range("E9").Value = '  +1.78%  '
$ws.Range("D10").Value = '0.146'
$ws.Range("E10").Value = '  +5.60%  '
$ws.Range("D11").Value = '7.32'
$ws.Range("E11").Value = '  +5.49%  '
$ws.Range("D12").Value = '0.437'
$ws.Range("E12").Value = '  +3.54%  '
$ws.Range("D13").Value = '0.0000221'
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = '33.09'
$ws.Range("E14").Value = '  +4.35%  '
$ws.Range("D15").Value = '4.145.65'
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '3.548.77'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = '68.097.74'
$ws.Range("E17").Value = '  +1.34%  '
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").Value = '6.76'
$ws.Range("E19").Value = '  +5.87%  '
$ws.Range("D20").Value = '15.90'
$ws.Range("E20").Value = '  +5.78%  '
$ws.Range("D21").Value = '9.97'
$ws.Range("E21").Value = '  +10.51%  '
$ws.Range("D22").Value = '454.63'
$ws.Range("E22").Value = '  +1.81%  '
$ws.Range("D23").Value = '0.640'
$ws.Range("E23").Value = '  +2.98%  '
$ws.Range("D24").Value = '78.20'
$ws.Range("E24").Value = '  +1.14%  '
$ws.Range("D25").Value = '10.52'
$ws.Range("E25").Value = '  +3.57%  '
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("D27").Value = '3.686.10'
$ws.Range("E27").Value = '  +1.50%  '
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("D29").Value = '9.07'
$ws.Range("E29").Value = '  +8.71%  '
$ws.Range("E30").Value = '  +3.40%  '
$ws.Range("E31").Value = '  +6.54%  '
$ws.Range("D32").Value = '0.170'
$ws.Range("E32").Value = '  +3.23%  '
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("D35").Value = '26.04'
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("E36").Value = '  +3.57%  '
$ws.Range("D37").Value = '3.539.00'
$ws.Range("E37").Value = '  +1.59%  '
$ws.Range("D38").Value = '8.25'
$ws.Range("E38").Value = '  +3.23%  '
$ws.Range("E39").Value = '  +7.35%  '
$ws.Range("D41").Value = '178.15'
$ws.Range("E41").Value = '  +3.93%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").Value = '0.0918'
$ws.Range("E42").Value = '  +5.47%  '
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("D44").Value = '5.58'
$ws.Range("E44").Value = '  +2.98%  '
$ws.Range("D45").Value = '30.79'
$ws.Range("E45").Value = '  +14.34%  '
$ws.Range("D46").Value = '0.894'
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("E47").Value = '  +6.67%  '
$ws.Range("D48").Value = '45.81'
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("E49").Value = '  +3.60%  '
$ws.Range("E50").Value = '  +3.30%  '
$ws.Range("E51").Value = '  +1.74%  '
